# The document has a first-page header (BTec logo) and two footers
# (default + first-page, both carrying the Pearson Edexcel logo). The
# commit simply swaps the internal "name" each embedded picture reports
# (wp:docPr / pic:cNvPr @name) between image1.* and image2.* - a pure
# rename, no visual/content change.
$d = $word.ActiveDocument
$sec = $d.Sections(1)

# --- Footers -------------------------------------------------------
# wdHeaderFooterPrimary (1)   -> footer2.xml -> image2.png -> image1.png
$footerPrimary = $sec.Footers(1)
if ($footerPrimary.Exists -and $footerPrimary.Range.InlineShapes.Count -ge 1) {
    $footerPrimary.Range.InlineShapes(1).Name = "image1.png"
}

# wdHeaderFooterFirstPage (2) -> footer1.xml -> image2.png -> image1.png
$footerFirst = $sec.Footers(2)
if ($footerFirst.Exists -and $footerFirst.Range.InlineShapes.Count -ge 1) {
    $footerFirst.Range.InlineShapes(1).Name = "image1.png"
}

# --- Header ----------------------------------------------------------
# wdHeaderFooterFirstPage (2) -> header1.xml -> image1.jpg -> image2.jpg
$headerFirst = $sec.Headers(2)
if ($headerFirst.Exists -and $headerFirst.Range.InlineShapes.Count -ge 1) {
    $headerFirst.Range.InlineShapes(1).Name = "image2.jpg"
}
